$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the date-formatted style (bold, centered, bordered, date number format)
# from the last existing row (343) down through the new rows (344:357),
# mirroring how the prior rows were formatted.
$ws.Range("A343").Copy()
$ws.Range("A344:A357").PasteSpecial(-4122)

# New daily case-count rows (10-23 Aug 2021)
$ws.Cells.Item(344, 1).Value = 44418
$ws.Cells.Item(344, 2).Value = 5
$ws.Cells.Item(344, 3).Value = 22
$ws.Cells.Item(344, 4).Value = 122.3717877405718

$ws.Cells.Item(345, 1).Value = 44419
$ws.Cells.Item(345, 2).Value = 2
$ws.Cells.Item(345, 3).Value = 21
$ws.Cells.Item(345, 4).Value = 116.809433752364

$ws.Cells.Item(346, 1).Value = 44420
$ws.Cells.Item(346, 2).Value = 6
$ws.Cells.Item(346, 3).Value = 27
$ws.Cells.Item(346, 4).Value = 150.1835576816109

$ws.Cells.Item(347, 1).Value = 44421
$ws.Cells.Item(347, 2).Value = 5
$ws.Cells.Item(347, 3).Value = 30
$ws.Cells.Item(347, 4).Value = 166.8706196462343

$ws.Cells.Item(348, 1).Value = 44422
$ws.Cells.Item(348, 2).Value = 2
$ws.Cells.Item(348, 3).Value = 30
$ws.Cells.Item(348, 4).Value = 166.8706196462343

$ws.Cells.Item(349, 1).Value = 44423
$ws.Cells.Item(349, 2).Value = 3
$ws.Cells.Item(349, 3).Value = 29
$ws.Cells.Item(349, 4).Value = 161.3082656580265

$ws.Cells.Item(350, 1).Value = 44424
$ws.Cells.Item(350, 2).Value = 5
$ws.Cells.Item(350, 3).Value = 28
$ws.Cells.Item(350, 4).Value = 155.7459116698187

$ws.Cells.Item(351, 1).Value = 44425
$ws.Cells.Item(351, 2).Value = 5
$ws.Cells.Item(351, 3).Value = 28
$ws.Cells.Item(351, 4).Value = 155.7459116698187

$ws.Cells.Item(352, 1).Value = 44426
$ws.Cells.Item(352, 2).Value = 0
$ws.Cells.Item(352, 3).Value = 26
$ws.Cells.Item(352, 4).Value = 144.621203693403

$ws.Cells.Item(353, 1).Value = 44427
$ws.Cells.Item(353, 2).Value = 7
$ws.Cells.Item(353, 3).Value = 27
$ws.Cells.Item(353, 4).Value = 150.1835576816109

$ws.Cells.Item(354, 1).Value = 44428
$ws.Cells.Item(354, 2).Value = 4
$ws.Cells.Item(354, 3).Value = 26
$ws.Cells.Item(354, 4).Value = 144.621203693403

$ws.Cells.Item(355, 1).Value = 44429
$ws.Cells.Item(355, 2).Value = 0
$ws.Cells.Item(355, 3).Value = 24
$ws.Cells.Item(355, 4).Value = 133.4964957169874

$ws.Cells.Item(356, 1).Value = 44430
$ws.Cells.Item(356, 2).Value = 0
$ws.Cells.Item(356, 3).Value = 21
$ws.Cells.Item(356, 4).Value = 116.809433752364

$ws.Cells.Item(357, 1).Value = 44431
$ws.Cells.Item(357, 2).Value = 1
$ws.Cells.Item(357, 3).Value = 17
$ws.Cells.Item(357, 4).Value = 94.56001779953276
